# Auto-generated script to apply the XLSX data refresh described by the diff.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for
# specific Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets,
# matching a scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value2 = 9143.333000000001
$ws.Range("J29").Value2 = 10248.375
$ws.Range("L29").Value2 = 30745.125
$ws.Range("N29").Value2 = -31307.125

$ws.Range("H31").Value2 = 1224.2
$ws.Range("I31").Value2 = 1280.25
$ws.Range("J31").Value2 = 1000
$ws.Range("K31").Value2 = 3840.75
$ws.Range("L31").Value2 = 3000
$ws.Range("M31").Value2 = -3610.75
$ws.Range("N31").Value2 = -3460

$ws.Range("H32").Value2 = 1868.1428
$ws.Range("I32").Value2 = 1932.6666
$ws.Range("J32").Value2 = 1819.75
$ws.Range("K32").Value2 = 1932.6666
$ws.Range("L32").Value2 = 1819.75
$ws.Range("M32").Value2 = -1606.6666
$ws.Range("N32").Value2 = -2471.75

$ws.Range("H34").Value2 = 3699
$ws.Range("I34").Value2 = 3699
$ws.Range("K34").Value2 = 3699
$ws.Range("M34").Value2 = -3496

$ws.Range("H36").Value2 = 3699
$ws.Range("I36").Value2 = 3699
$ws.Range("K36").Value2 = 3699
$ws.Range("M36").Value2 = -2984

$ws.Range("H40").Value2 = 4769.9
$ws.Range("I40").Value2 = 3300
$ws.Range("J40").Value2 = 5399.857
$ws.Range("K40").Value2 = 3300
$ws.Range("L40").Value2 = 5399.857
$ws.Range("M40").Value2 = -3125
$ws.Range("N40").Value2 = -5749.857

$ws.Range("H55").Value2 = 2482.2
$ws.Range("I55").Value2 = 486.8
$ws.Range("K55").Value2 = 486.8
$ws.Range("M55").Value2 = -272.8

$ws.Range("H74").Value2 = 3862.625
$ws.Range("I74").Value2 = 3862.625
$ws.Range("K74").Value2 = 3862.625
$ws.Range("M74").Value2 = -2926.625

$ws.Range("H77").Value2 = 3862.625
$ws.Range("I77").Value2 = 3862.625
$ws.Range("K77").Value2 = 19313.125
$ws.Range("M77").Value2 = -14633.125

$ws.Range("H127").Value2 = 5731.9165
$ws.Range("I127").Value2 = 3222.875
$ws.Range("K127").Value2 = 9668.625
$ws.Range("M127").Value2 = -4708.625

$ws.Range("H137").Value2 = 4064.861
$ws.Range("I137").Value2 = 2636.6
$ws.Range("K137").Value2 = 7909.799999999999
$ws.Range("M137").Value2 = -5359.799999999999

$ws.Range("H138").Value2 = 3334.5715
$ws.Range("J138").Value2 = 3477.4119
$ws.Range("L138").Value2 = 10432.2357
$ws.Range("N138").Value2 = -20712.2357


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 7814367
$ws.Range("I32").Value2 = 8334491
$ws.Range("K32").Value2 = 8334491
$ws.Range("M32").Value2 = -8334204

$ws.Range("H50").Value2 = 1248.5
$ws.Range("I50").Value2 = 1199.3334
$ws.Range("J50").Value2 = 1322.25
$ws.Range("K50").Value2 = 1199.3334
$ws.Range("L50").Value2 = 1322.25
$ws.Range("M50").Value2 = -485.3334
$ws.Range("N50").Value2 = -2750.25

$ws.Range("H74").Value2 = 10006950
$ws.Range("I74").Value2 = 16670753
$ws.Range("J74").Value2 = 919945.6
$ws.Range("K74").Value2 = 16670753
$ws.Range("L74").Value2 = 919945.6
$ws.Range("M74").Value2 = -16669879
$ws.Range("N74").Value2 = -921693.6

$ws.Range("H77").Value2 = 10006950
$ws.Range("I77").Value2 = 16670753
$ws.Range("J77").Value2 = 919945.6
$ws.Range("K77").Value2 = 83353765
$ws.Range("L77").Value2 = 4599728
$ws.Range("M77").Value2 = -83349397
$ws.Range("N77").Value2 = -4608464

$ws.Range("H102").Value2 = 7077.8823
$ws.Range("I102").Value2 = 6770.25
$ws.Range("J102").Value2 = 12000
$ws.Range("K102").Value2 = 6770.25
$ws.Range("L102").Value2 = 12000
$ws.Range("M102").Value2 = -5148.25
$ws.Range("N102").Value2 = -15244

$ws.Range("H122").Value2 = 3150.2222
$ws.Range("I122").Value2 = 2724.4614
$ws.Range("K122").Value2 = 8173.3842
$ws.Range("M122").Value2 = -5723.3842

$ws.Range("H132").Value2 = 3995.9355
$ws.Range("I132").Value2 = 1803.1305
$ws.Range("J132").Value2 = 10300.25
$ws.Range("K132").Value2 = 5409.3915
$ws.Range("L132").Value2 = 30900.75
$ws.Range("M132").Value2 = -2879.3915
$ws.Range("N132").Value2 = -35960.75


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value2 = 97436.5
$ws.Range("I40").Value2 = 97429
$ws.Range("K40").Value2 = 97429
$ws.Range("M40").Value2 = -97164

$ws.Range("H134").Value2 = 386154.16
$ws.Range("I134").Value2 = 1304.1305
$ws.Range("K134").Value2 = 3912.3915
$ws.Range("M134").Value2 = -1377.3915


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 592421.4
$ws.Range("I31").Value2 = 18464.916
$ws.Range("K31").Value2 = 18464.916
$ws.Range("M31").Value2 = -18169.916

$ws.Range("H34").Value2 = 592421.4
$ws.Range("I34").Value2 = 18464.916
$ws.Range("K34").Value2 = 18464.916
$ws.Range("M34").Value2 = -18262.916

$ws.Range("H127").Value2 = 72000
$ws.Range("I127").Value2 = 0
$ws.Range("J127").Value2 = 72000
$ws.Range("K127").Value2 = 0
$ws.Range("M127").Value2 = 72000
$ws.Range("N127").Value2 = -81920
$ws.Range("L127").ClearContents()

$ws.Range("H132").Value2 = 2049.4707
$ws.Range("I132").Value2 = 2162.875
$ws.Range("J132").Value2 = 235
$ws.Range("K132").Value2 = 6488.625
$ws.Range("L132").Value2 = 705
$ws.Range("M132").Value2 = -3958.625
$ws.Range("N132").Value2 = -5765

$ws.Range("H134").Value2 = 287197.47
$ws.Range("I134").Value2 = 314011.4
$ws.Range("J134").Value2 = 1182
$ws.Range("K134").Value2 = 942034.2000000001
$ws.Range("L134").Value2 = 3546
$ws.Range("M134").Value2 = -939499.2000000001
$ws.Range("N134").Value2 = -8616

$ws.Range("H141").Value2 = 181237.25
$ws.Range("J141").Value2 = 181237.25
$ws.Range("L141").Value2 = 181237.25
$ws.Range("N141").Value2 = -191597.25


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value2 = 118.083336
$ws.Range("I17").Value2 = 51.285713
$ws.Range("J17").Value2 = 211.6
$ws.Range("K17").Value2 = 153.857139
$ws.Range("L17").Value2 = 634.8
$ws.Range("M17").Value2 = 15.14286099999998
$ws.Range("N17").Value2 = -972.8

$ws.Range("H140").Value2 = 160034.31
$ws.Range("I140").Value2 = 160034.31
$ws.Range("K140").Value2 = 480102.93
$ws.Range("M140").Value2 = -474922.93


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value2 = 100833.336
$ws.Range("I35").Value2 = 116666.664
$ws.Range("K35").Value2 = 116666.664
$ws.Range("M35").Value2 = -116368.664

$ws.Range("H36").Value2 = 0
$ws.Range("I36").Value2 = 0
$ws.Range("J36").Value2 = 0
$ws.Range("K36").Value2 = 0
$ws.Range("N36").Value2 = 0
$ws.Range("L36").ClearContents()
$ws.Range("M36").ClearContents()

$ws.Range("H75").Value2 = 0
$ws.Range("J75").Value2 = 0
$ws.Range("N75").Value2 = 0
$ws.Range("L75").ClearContents()

$ws.Range("H78").Value2 = 0
$ws.Range("J78").Value2 = 0
$ws.Range("N78").Value2 = 0
$ws.Range("L78").ClearContents()

$ws.Range("H102").Value2 = 3162.5454
$ws.Range("I102").Value2 = 4131.1665
$ws.Range("J102").Value2 = 2000.2
$ws.Range("K102").Value2 = 4131.1665
$ws.Range("L102").Value2 = 2000.2
$ws.Range("M102").Value2 = -2509.1665
$ws.Range("N102").Value2 = -5244.2

$ws.Range("H126").Value2 = 4997.2856
$ws.Range("I126").Value2 = 4800.1113
$ws.Range("K126").Value2 = 14400.3339
$ws.Range("M126").Value2 = -11930.3339

$ws.Range("H132").Value2 = 43486490
$ws.Range("I132").Value2 = 50002960
$ws.Range("K132").Value2 = 150008880
$ws.Range("M132").Value2 = -150006350

$ws.Range("H133").Value2 = 0
$ws.Range("J133").Value2 = 0
$ws.Range("N133").Value2 = 0
$ws.Range("L133").ClearContents()

$ws.Range("H136").Value2 = 28970.715
$ws.Range("J136").Value2 = 28970.715
$ws.Range("L136").Value2 = 86912.145
$ws.Range("N136").Value2 = -92012.145


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value2 = 1905.8
$ws.Range("I35").Value2 = 1905.8
$ws.Range("J35").Value2 = 0
$ws.Range("K35").Value2 = 1905.8
$ws.Range("L35").Value2 = 0
$ws.Range("N35").Value2 = -1569.8
$ws.Range("M35").ClearContents()

$ws.Range("H61").Value2 = 1479.3572
$ws.Range("I61").Value2 = 1542.0769
$ws.Range("J61").Value2 = 664
$ws.Range("K61").Value2 = 1542.0769
$ws.Range("L61").Value2 = 664
$ws.Range("M61").Value2 = -1340.0769
$ws.Range("N61").Value2 = -1068

$ws.Range("H100").Value2 = 3158.3157
$ws.Range("I100").Value2 = 2667.3333
$ws.Range("K100").Value2 = 2667.3333
$ws.Range("M100").Value2 = -2126.3333

$ws.Range("H113").Value2 = 1479.3572
$ws.Range("I113").Value2 = 1542.0769
$ws.Range("J113").Value2 = 664
$ws.Range("K113").Value2 = 1542.0769
$ws.Range("L113").Value2 = 664
$ws.Range("M113").Value2 = 627.9231
$ws.Range("N113").Value2 = -5004


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value2 = 0
$ws.Range("I29").Value2 = 0
$ws.Range("K29").Value2 = 0
$ws.Range("M29").ClearContents()

$ws.Range("H31").Value2 = 0
$ws.Range("I31").Value2 = 0
$ws.Range("J31").Value2 = 0
$ws.Range("K31").Value2 = 0
$ws.Range("N31").Value2 = 0
$ws.Range("L31").ClearContents()
$ws.Range("M31").ClearContents()

$ws.Range("H75").Value2 = 11950028
$ws.Range("J75").Value2 = 10800032
$ws.Range("L75").Value2 = 10800032
$ws.Range("N75").Value2 = -10801904

$ws.Range("H78").Value2 = 11950028
$ws.Range("J78").Value2 = 10800032
$ws.Range("L78").Value2 = 32400096
$ws.Range("N78").Value2 = -32409456

$ws.Range("H132").Value2 = 2440.3333
$ws.Range("I132").Value2 = 2026.2963
$ws.Range("K132").Value2 = 6078.8889
$ws.Range("M132").Value2 = -3548.8889

